# Commit: Add SQQ trade, add more info in buy leg results
#
# This script replaces the old "buy UNG put / write UNG put" trade sheet
# with a new "buy SQQQ call" trade, and removes the now-stale "write" leg
# values (rows 19-30 keep their FIELD labels in column A, but column B is
# cleared out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab (and the <sheet> entry in workbook.xml)
$ws.Name = "20160928 UNG"

# Move the active selection to F25, as recorded in the saved view state
$ws.Range("F25").Select()

# --- Buy leg (rows 2-18): new SQQQ call trade -----------------------------
$ws.Cells.Item(2, 2).Value2 = "20160928 +SQQQ-161021C13.00"   # option_buy
$ws.Cells.Item(4, 2).Value2 = "call"                            # option_type
$ws.Cells.Item(5, 2).Value2 = 13                                # strike_buy
$ws.Cells.Item(9, 2).Value2 = 0.85                              # premium_buy
$ws.Cells.Item(10, 2).Value2 = 13.45                            # underlying_buy
$ws.Cells.Item(14, 2).Value2 = 13                               # entry_date_buy_HH
$ws.Cells.Item(15, 2).Value2 = 38                               # entry_date_buy_MM
$ws.Cells.Item(16, 2).Value2 = 19                               # entry_date_buy_SS
$ws.Cells.Item(17, 2).Value2 = "SQQQ"                           # ticker
$ws.Cells.Item(18, 2).Value2 = 0.4314                           # historical_volatility_buy

# --- Write leg (rows 19-30): clear stale values, keep only the labels ----
$ws.Range("B19:B30").ClearContents()
